# "generating data on the fly"
# Fill in the second block (columns H:L) of the two "conv" mini-tables
# (rows 15-20 and rows 23-28) on Sheet1, mirroring the existing
# B:F / H:M / P:T / V:AA pattern already used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlCenter = -4108

# ---------------------------------------------------------------------
# Block 1: rows 15-20, columns H:L
# ---------------------------------------------------------------------

# Make sure the whole block carries the centred style used by the rest
# of the sheet's mini-tables (matches cellXfs index 1).
$ws.Range("H15:L20").HorizontalAlignment = $xlCenter

# Header row
$ws.Range("H15").Value = "in"
$ws.Range("I15").Value = "kernel"
$ws.Range("J15").Value = "padding"
$ws.Range("K15").Value = "stride"
$ws.Range("L15").Value = "out"

# Row 16 - first data row: literal input size, then kernel/padding/stride
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 2
$ws.Range("L16").Formula = "=ROUNDDOWN((H16+2*J16-(I16-1)-1)/K16+1,0)"

# Row 17 - input chains off previous row's output
$ws.Range("H17").Formula = "=L16"
$ws.Range("I17").Value = 5
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 2
$ws.Range("L17").Formula = "=ROUNDDOWN((H17+2*J17-(I17-1)-1)/K17+1,0)"

# Row 18
$ws.Range("H18").Formula = "=L17"
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 2
$ws.Range("L18").Formula = "=ROUNDDOWN((H18+2*J18-(I18-1)-1)/K18+1,0)"

# Row 19
$ws.Range("H19").Formula = "=L18"
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 2
$ws.Range("L19").Formula = "=ROUNDDOWN((H19+2*J19-(I19-1)-1)/K19+1,0)"

# Row 20
$ws.Range("H20").Formula = "=L19"
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 2
$ws.Range("L20").Formula = "=ROUNDDOWN((H20+2*J20-(I20-1)-1)/K20+1,0)"

# Row 20's height changes 15 -> 13.8 to match the data rows above it
$ws.Rows.Item(20).RowHeight = 13.8

# ---------------------------------------------------------------------
# Block 2: rows 23-28, columns H:L
# ---------------------------------------------------------------------

$ws.Range("H23:L28").HorizontalAlignment = $xlCenter

# Header row
$ws.Range("H23").Value = "in"
$ws.Range("I23").Value = "kernel"
$ws.Range("J23").Value = "padding"
$ws.Range("K23").Value = "stride"
$ws.Range("L23").Value = "out"

# Row 24 - first data row of second block
$ws.Range("H24").Value = 229
$ws.Range("I24").Value = 5
$ws.Range("J24").Value = 1
$ws.Range("K24").Value = 2
$ws.Range("L24").Formula = "=ROUNDDOWN((H24+2*J24-(I24-1)-1)/K24+1,0)"

# Row 25
$ws.Range("H25").Formula = "=L24"
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = 2
$ws.Range("L25").Formula = "=ROUNDDOWN((H25+2*J25-(I25-1)-1)/K25+1,0)"

# Row 26
$ws.Range("H26").Formula = "=L25"
$ws.Range("I26").Value = 5
$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 2
$ws.Range("L26").Formula = "=ROUNDDOWN((H26+2*J26-(I26-1)-1)/K26+1,0)"

# Row 27
$ws.Range("H27").Formula = "=L26"
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 2
$ws.Range("L27").Formula = "=ROUNDDOWN((H27+2*J27-(I27-1)-1)/K27+1,0)"

# Row 28 - note: padding/stride go back to 2/2 here (matches diff)
$ws.Range("H28").Formula = "=L27"
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 2
$ws.Range("L28").Formula = "=ROUNDDOWN((H28+2*J28-(I28-1)-1)/K28+1,0)"

# Row heights 27/28 change 15 -> 13.8 to match the data rows above them
$ws.Rows.Item(27).RowHeight = 13.8
$ws.Rows.Item(28).RowHeight = 13.8

# ---------------------------------------------------------------------
# View state: scroll position + active selection moved
# ---------------------------------------------------------------------
[void]$ws.Range("A12").Select()
[void]$ws.Range("K29").Select()
